# Carga Instalada.xlsx - update service entrance / load calculations
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Updated nominal power (kW) readings for the motor loads, recalculated
# --- using the WEG efficiency/power-factor table (per commit message).
$ws.Range("D6").Value  = 2.76    # Misturador Submerso
$ws.Range("D7").Value  = 0.9     # Bomba Helicoidal
$ws.Range("D8").Value  = 3.68    # Bomba Centrifuga Submersivel
$ws.Range("D9").Value  = 10.77   # Soprador
$ws.Range("D10").Value = 0.3     # Bomba Dosadora com Diafragma

# F column holds =E*D formulas and the F11 grand total is a SUM formula, so
# they recalculate automatically from the new D values.

# --- Emphasize the installed-load total row: bold the (otherwise empty)
# --- D11 cell and give it the same "0.00 kW" number format used elsewhere.
$ws.Range("D11").Font.Bold = $true
$ws.Range("D11").NumberFormat = '0.00\ "kW"'

# --- Restore the working selection to the data table.
$ws.Range("A2:F11").Select() | Out-Null
